$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the BMP variation coefficient value (row 5, column B)
$ws.Range("B5").Value = 14.51

# Update the reference text for the BMP row to point to the new source
$ws.Range("C5").Value = "Hafner (2020), Tab. 3 (mean of all 4 substrates SA-SD without cellulose)"

# Move the active selection to B6, matching the saved workbook state
$ws.Range("B6").Select()
